# Auto-generated script applying the cryptos price/volume update
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D (Price) updates
$priceUpdates = @{
    "D2" = "63.863.92"
    "D3" = "2.629.29"
    "D5" = "598.16"
    "D6" = "150.58"
    "D8" = "0.589"
    "D10" = "5.71"
    "D11" = "0.383"
    "D12" = "0.151"
    "D13" = "27.73"
    "D14" = "3.097.99"
    "D15" = "63.674.49"
    "D16" = "0.0000150"
    "D17" = "2.632.06"
    "D18" = "12.33"
    "D19" = "4.64"
    "D20" = "349.50"
    "D21" = "6.89"
    "D23" = "5.71"
    "D26" = "9.21"
    "D27" = "1.68"
    "D28" = "567.98"
    "D29" = "8.28"
    "D33" = "0.0₃0848"
    "D34" = "1.74"
    "D35" = "5.25"
    "D36" = "169.35"
    "D37" = "0.409"
    "D38" = "0.999"
    "D39" = "1.95"
    "D40" = "19.38"
    "D42" = "170.37"
    "D43" = "39.90"
    "D45" = "0.0600"
    "D46" = "21.45"
    "D48" = "0.0248"
    "D49" = "1.99"
    "D50" = "0.0970"
    "D51" = "19.21"
}

# Column D cells whose new value looks numeric need to be forced to
# Text so Excel does not reinterpret them (loses the distinct period-
# grouped / trailing-zero text formatting used in the source data).
$priceForceText = @(
    "D5"
    "D6"
    "D8"
    "D10"
    "D11"
    "D12"
    "D13"
    "D16"
    "D18"
    "D19"
    "D20"
    "D21"
    "D23"
    "D26"
    "D27"
    "D28"
    "D29"
    "D34"
    "D35"
    "D36"
    "D37"
    "D38"
    "D39"
    "D40"
    "D42"
    "D43"
    "D45"
    "D46"
    "D48"
    "D49"
    "D50"
    "D51"
)

foreach ($cell in $priceUpdates.Keys) {
    if ($priceForceText -contains $cell) {
        $ws.Range($cell).NumberFormat = "@"
    }
    $ws.Range($cell).Value = $priceUpdates[$cell]
}

# Column E (Volume/1h change %) updates - these are already padded
# "  +x.xx%  " style text strings, never numeric-parseable.
$volumeUpdates = @{
    "E2" = "  +0.23%  "
    "E3" = "  -0.99%  "
    "E4" = "  -0.04%  "
    "E5" = "  -0.87%  "
    "E6" = "  +1.85%  "
    "E7" = "  -0.04%  "
    "E8" = "  +0.14%  "
    "E9" = "  +0.66%  "
    "E10" = "  +2.01%  "
    "E11" = "  +3.44%  "
    "E12" = "  -1.10%  "
    "E13" = "  +0.38%  "
    "E14" = "  -1.18%  "
    "E15" = "  +0.01%  "
    "E16" = "  +2.53%  "
    "E17" = "  -1.16%  "
    "E18" = "  +7.56%  "
    "E19" = "  +1.70%  "
    "E20" = "  +1.95%  "
    "E21" = "  -1.54%  "
    "E22" = "  -0.30%  "
    "E23" = "  +2.28%  "
    "E24" = "  -0.89%  "
    "E25" = "  +12.61%  "
    "E26" = "  +1.24%  "
    "E27" = "  -0.74%  "
    "E28" = "  +1.96%  "
    "E29" = "  +4.28%  "
    "E30" = "  +0.29%  "
    "E31" = "  +0.26%  "
    "E32" = "  +1.13%  "
    "E33" = "  +3.20%  "
    "E34" = "  -0.88%  "
    "E35" = "  +0.88%  "
    "E36" = "  +1.37%  "
    "E37" = "  +0.68%  "
    "E38" = "  -0.09%  "
    "E39" = "  +0.67%  "
    "E40" = "  +1.10%  "
    "E41" = "  +0.04%  "
    "E42" = "  +1.31%  "
    "E43" = "  -0.22%  "
    "E44" = "  +3.22%  "
    "E45" = "  +4.11%  "
    "E46" = "  -3.78%  "
    "E47" = "  -0.16%  "
    "E48" = "  +0.01%  "
    "E49" = "  +6.54%  "
    "E50" = "  +0.57%  "
    "E51" = "  +2.05%  "
}

foreach ($cell in $volumeUpdates.Keys) {
    $ws.Range($cell).Value = $volumeUpdates[$cell]
}
